$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Correct the NHC value in row 3 (C3): 1 -> 2
$ws.Range("C3").Value = 2

# Move/leave the active selection on C3 (matches the saved view state)
$ws.Range("C3").Select()
